$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 909.8570999999999
$ws.Range("I2").Value = 314
$ws.Range("K2").Value = 314
$ws.Range("M2").Value = -201
$ws.Range("H28").Value = 2236.111
$ws.Range("I28").Value = 1265.625
$ws.Range("K28").Value = 1265.625
$ws.Range("M28").Value = -780.625
$ws.Range("H51").Value = 12399.934
$ws.Range("J51").Value = 12399.934
$ws.Range("L51").Value = 12399.934
$ws.Range("N51").Value = -13367.934
$ws.Range("H80").Value = 886.9524
$ws.Range("I80").Value = 757.1667
$ws.Range("J80").Value = 1060
$ws.Range("K80").Value = 2271.5001
$ws.Range("L80").Value = 3180
$ws.Range("M80").Value = -1273.5001
$ws.Range("N80").Value = -5176
$ws.Range("H83").Value = 886.9524
$ws.Range("I83").Value = 757.1667
$ws.Range("J83").Value = 1060
$ws.Range("K83").Value = 6814.5003
$ws.Range("L83").Value = 9540
$ws.Range("M83").Value = -1822.5003
$ws.Range("N83").Value = -19524
$ws.Range("H106").Value = 3995.8235
$ws.Range("I106").Value = 3295.3333
$ws.Range("K106").Value = 3295.3333
$ws.Range("M106").Value = -2664.3333
$ws.Range("H129").Value = 2949.5
$ws.Range("I129").Value = 2266.3333
$ws.Range("K129").Value = 6798.999899999999
$ws.Range("M129").Value = -1798.999899999999
$ws.Range("H131").Value = 1925
$ws.Range("I131").Value = 1925
$ws.Range("K131").Value = 5775
$ws.Range("M131").Value = -735
$ws.Range("H141").Value = 3335.28
$ws.Range("I141").Value = 1887.5
$ws.Range("K141").Value = 5662.5
$ws.Range("M141").Value = -482.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 7999.5
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 7999.5
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 7999.5
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -8339.5
$ws.Range("H32").Value = 6219.076
$ws.Range("I32").Value = 2891.41
$ws.Range("J32").Value = 17496.166
$ws.Range("K32").Value = 2891.41
$ws.Range("L32").Value = 17496.166
$ws.Range("M32").Value = -2604.41
$ws.Range("N32").Value = -18070.166
$ws.Range("H46").Value = 4000
$ws.Range("J46").Value = 4000
$ws.Range("L46").Value = 4000
$ws.Range("N46").Value = -4638
$ws.Range("H61").Value = 4355
$ws.Range("I61").Value = 3740.5
$ws.Range("K61").Value = 3740.5
$ws.Range("M61").Value = -3528.5
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("K123").Value = 0
$ws.Range("M123").ClearContents()
$ws.Range("H136").Value = 4355
$ws.Range("I136").Value = 3740.5
$ws.Range("K136").Value = 11221.5
$ws.Range("M136").Value = -8671.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H98").Value = 125000
$ws.Range("J98").Value = 125000
$ws.Range("L98").Value = 125000
$ws.Range("N98").Value = -130990
$ws.Range("H99").Value = 1793.6364
$ws.Range("I99").Value = 1153.8
$ws.Range("K99").Value = 1153.8
$ws.Range("M99").Value = 344.2
$ws.Range("H134").Value = 2600.3572
$ws.Range("I134").Value = 2600.3572
$ws.Range("K134").Value = 7801.071599999999
$ws.Range("M134").Value = -5266.071599999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13086.667
$ws.Range("I31").Value = 4809.222
$ws.Range("K31").Value = 4809.222
$ws.Range("M31").Value = -4514.222
$ws.Range("H34").Value = 13086.667
$ws.Range("I34").Value = 4809.222
$ws.Range("K34").Value = 4809.222
$ws.Range("M34").Value = -4607.222
$ws.Range("H105").Value = 4001617
$ws.Range("I105").Value = 4001617
$ws.Range("K105").Value = 4001617
$ws.Range("M105").Value = -3999870
$ws.Range("H132").Value = 31515.117
$ws.Range("I132").Value = 33359.875
$ws.Range("K132").Value = 100079.625
$ws.Range("M132").Value = -97549.625
$ws.Range("H141").Value = 322022.22
$ws.Range("J141").Value = 555640
$ws.Range("L141").Value = 555640
$ws.Range("N141").Value = -566000

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 1519.5
$ws.Range("I29").Value = 3430
$ws.Range("K29").Value = 10290
$ws.Range("M29").Value = -10013
$ws.Range("H48").Value = 5999
$ws.Range("I48").Value = 5999
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 17997
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -17747
$ws.Range("N48").ClearContents()
$ws.Range("H129").Value = 1033.7778
$ws.Range("I129").Value = 1033.7778
$ws.Range("K129").Value = 3101.3334
$ws.Range("M129").Value = 1898.6666
$ws.Range("H131").Value = 1576.0416
$ws.Range("I131").Value = 1585
$ws.Range("J131").Value = 1572.3529
$ws.Range("K131").Value = 4755
$ws.Range("L131").Value = 4717.0587
$ws.Range("M131").Value = 285
$ws.Range("N131").Value = -14797.0587
$ws.Range("H138").Value = 2797.6
$ws.Range("I138").Value = 2797.6
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 8392.799999999999
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -3252.799999999999
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 3259.4
$ws.Range("I139").Value = 1382.5
$ws.Range("J139").Value = 6074.75
$ws.Range("K139").Value = 4147.5
$ws.Range("L139").Value = 18224.25
$ws.Range("M139").Value = 992.5
$ws.Range("N139").Value = -28504.25
$ws.Range("H140").Value = 2028.2916
$ws.Range("I140").Value = 1511.25
$ws.Range("K140").Value = 4533.75
$ws.Range("M140").Value = 646.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 9999
$ws.Range("J27").Value = 9999
$ws.Range("L27").Value = 9999
$ws.Range("N27").Value = -10331
$ws.Range("H132").Value = 4945
$ws.Range("I132").Value = 4973.2144
$ws.Range("K132").Value = 14919.6432
$ws.Range("M132").Value = -12389.6432
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5352.6665
$ws.Range("I61").Value = 4810.8945
$ws.Range("K61").Value = 4810.8945
$ws.Range("M61").Value = -4608.8945
$ws.Range("H113").Value = 5352.6665
$ws.Range("I113").Value = 4810.8945
$ws.Range("K113").Value = 4810.8945
$ws.Range("M113").Value = -2640.8945
$ws.Range("H132").Value = 153368340
$ws.Range("I132").Value = 153368340
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 460105020
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -460102490
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2999
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H100").Value = 1693.65
$ws.Range("I100").Value = 1727.4375
$ws.Range("K100").Value = 3454.875
$ws.Range("M100").Value = -2913.875
